# Update "想去人数" (column F) figures on both the "展览" and "全部类型"
# sheets. These two sheets hold duplicate data, so the same F-column
# cells are updated identically on each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    4  = 1604
    5  = 614
    8  = 11425
    9  = 25
    12 = 353
    13 = 1089
    14 = 792
    15 = 12362
    16 = 13029
    21 = 15
    22 = 86
    24 = 102
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
